$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - zh-cn handback finished in-place (handoff file == handback file)
#   - de-de handback finished in-place too, with its own handback timestamp
#   - Both locale sheets gain "Latest Target File" / "Latest Handback File"
#     columns (F/G) pointing at the same files that were handed off
# ---------------------------------------------------------------------------

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: refresh the status column for both rows ------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusText
$ov.Range("C2").Value = $statusText
$ov.Range("B3").Value = $statusText
$ov.Range("C3").Value = $statusText

# ---- zh-cn sheet ----------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhA2md  = "https://github.com/OpenLocalizationTest/oltest/blob/d46a6ce8bc404bc81a022195c080a79b3db01003/e2e/398da73c-b7db-4880-9842-02f7a0f208e3.md"
$zhD2xlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ac8f9fa6683d2000383657573993399136320e6/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/398da73c-b7db-4880-9842-02f7a0f208e3.e317bd12f3c10c39323b182a24439411041eb136.zh-cn.xlf"
$zhA3md  = "https://github.com/OpenLocalizationTest/oltest/blob/d46a6ce8bc404bc81a022195c080a79b3db01003/e2e/c21a49a6-aec9-42aa-b340-a6d52e78f9db.md"
$zhD3xlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ac8f9fa6683d2000383657573993399136320e6/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/c21a49a6-aec9-42aa-b340-a6d52e78f9db.bb6593a1a67dfff0c6fd931b7bc9625e912a488c.zh-cn.xlf"

$zhA2disp = "398da73c-b7db-4880-9842-02f7a0f208e3.md"
$zhD2disp = "398da73c-b7db-4880-9842-02f7a0f208e3.e317bd12f3c10c39323b182a24439411041eb136.zh-cn.xlf"
$zhA3disp = "c21a49a6-aec9-42aa-b340-a6d52e78f9db.md"
$zhD3disp = "c21a49a6-aec9-42aa-b340-a6d52e78f9db.bb6593a1a67dfff0c6fd931b7bc9625e912a488c.zh-cn.xlf"

# status + handback datetime text (shared across both rows)
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText
$zh.Range("H2").Value = "2016-03-30 10:49:55"
$zh.Range("H3").Value = "2016-03-30 10:49:55"

# Latest Target File / Latest Handback File columns - the handback file is
# the same xlf that was last handed off, since the locale is in sync.
$zh.Range("F2").Value = $zhA2disp
$zh.Range("G2").Value = $zhD2disp
$zh.Range("F3").Value = $zhA3disp
$zh.Range("G3").Value = $zhD3disp

# Rebuild the hyperlinks so the new F/G links land between the row-2 and
# row-3 links (matches the order Excel lists them in after regenerating).
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhA2md,  "", "", $zhA2disp)
$zh.Hyperlinks.Add($zh.Range("D2"), $zhD2xlf, "", "", $zhD2disp)
$zh.Hyperlinks.Add($zh.Range("F2"), $zhA2md,  "", "", $zhA2disp)
$zh.Hyperlinks.Add($zh.Range("G2"), $zhD2xlf, "", "", $zhD2disp)
$zh.Hyperlinks.Add($zh.Range("A3"), $zhA3md,  "", "", $zhA3disp)
$zh.Hyperlinks.Add($zh.Range("D3"), $zhD3xlf, "", "", $zhD3disp)
$zh.Hyperlinks.Add($zh.Range("F3"), $zhA3md,  "", "", $zhA3disp)
$zh.Hyperlinks.Add($zh.Range("G3"), $zhD3xlf, "", "", $zhD3disp)

# ---- de-de sheet ------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$deA2md  = "https://github.com/OpenLocalizationTest/oltest/blob/d46a6ce8bc404bc81a022195c080a79b3db01003/e2e/398da73c-b7db-4880-9842-02f7a0f208e3.md"
$deD2xlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66074f2c392ed15f164173102cf40cdebdcc2cf2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/398da73c-b7db-4880-9842-02f7a0f208e3.e317bd12f3c10c39323b182a24439411041eb136.de-de.xlf"
$deA3md  = "https://github.com/OpenLocalizationTest/oltest/blob/d46a6ce8bc404bc81a022195c080a79b3db01003/e2e/c21a49a6-aec9-42aa-b340-a6d52e78f9db.md"
$deD3xlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66074f2c392ed15f164173102cf40cdebdcc2cf2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/c21a49a6-aec9-42aa-b340-a6d52e78f9db.bb6593a1a67dfff0c6fd931b7bc9625e912a488c.de-de.xlf"

$deA2disp = "398da73c-b7db-4880-9842-02f7a0f208e3.md"
$deD2disp = "398da73c-b7db-4880-9842-02f7a0f208e3.e317bd12f3c10c39323b182a24439411041eb136.de-de.xlf"
$deA3disp = "c21a49a6-aec9-42aa-b340-a6d52e78f9db.md"
$deD3disp = "c21a49a6-aec9-42aa-b340-a6d52e78f9db.bb6593a1a67dfff0c6fd931b7bc9625e912a488c.de-de.xlf"

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText
$de.Range("H2").Value = "2016-03-30 10:50:13"
$de.Range("H3").Value = "2016-03-30 10:50:13"

$de.Range("F2").Value = $deA2disp
$de.Range("G2").Value = $deD2disp
$de.Range("F3").Value = $deA3disp
$de.Range("G3").Value = $deD3disp

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deA2md,  "", "", $deA2disp)
$de.Hyperlinks.Add($de.Range("D2"), $deD2xlf, "", "", $deD2disp)
$de.Hyperlinks.Add($de.Range("F2"), $deA2md,  "", "", $deA2disp)
$de.Hyperlinks.Add($de.Range("G2"), $deD2xlf, "", "", $deD2disp)
$de.Hyperlinks.Add($de.Range("A3"), $deA3md,  "", "", $deA3disp)
$de.Hyperlinks.Add($de.Range("D3"), $deD3xlf, "", "", $deD3disp)
$de.Hyperlinks.Add($de.Range("F3"), $deA3md,  "", "", $deA3disp)
$de.Hyperlinks.Add($de.Range("G3"), $deD3xlf, "", "", $deD3disp)
